$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the emoji/color marker strings used in column A (statut) and
# the matching label "noir" -> "bleu" used in column B (statut_label),
# mirroring the shared-strings table edit:
#   🟥 -> 📕 (rouge)
#   ⬛ -> 📘 (noir -> bleu)
#   🟧 -> 📙 (orange)

$usedRange = $ws.UsedRange

$usedRange.Replace("🟥", "📕", 1) | Out-Null
$usedRange.Replace("⬛", "📘", 1) | Out-Null
$usedRange.Replace("🟧", "📙", 1) | Out-Null
$usedRange.Replace("noir", "bleu", 1) | Out-Null
